# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# Reorders the "Periodo Mora" / "Valor Mora" rows (B16:J23 block) so the
# periods run oldest -> newest (1912, 2001, 2002, 2003, 2011, 2012, 2101,
# 2102) instead of newest -> oldest, carrying each row's "Valor Mora"
# along with its matching period.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Periodo Mora" (column E) values, row 16 .. row 23, oldest first.
$periodos = @("1912", "2001", "2002", "2003", "2011", "2012", "2101", "2102")

# New "Valor Mora" (column F) values, row 16 .. row 23, matching the
# periods above.
$valores = @(68000, 68000, 68000, 68000, 35112, 35112, 35112, 35112)

for ($i = 0; $i -lt 8; $i++) {
    $row = 16 + $i
    $ws.Cells.Item($row, 5).Value = $periodos[$i]
    $ws.Cells.Item($row, 6).Value = $valores[$i]
}
